$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with refreshed crypto data.
# Leading apostrophe forces text interpretation (values like '1.002' or
# dotted-thousands prices must stay literal text, matching the source data);
# Style is reset to Normal afterwards so no stray number formatting sticks.
$ws.Range('D2').Value = "'24.530.05"
$ws.Range('D2').Style = "Normal"
$ws.Range('E2').Value = "'  +3.18%  "
$ws.Range('E2').Style = "Normal"
$ws.Range('D3').Value = "'1.692.61"
$ws.Range('D3').Style = "Normal"
$ws.Range('E3').Value = "'  +1.57%  "
$ws.Range('E3').Style = "Normal"
$ws.Range('D4').Value = "'1.002"
$ws.Range('D4').Style = "Normal"
$ws.Range('E4').Value = "'  +0.30%  "
$ws.Range('E4').Style = "Normal"
$ws.Range('D5').Value = "'314.29"
$ws.Range('D5').Style = "Normal"
$ws.Range('E5').Value = "'  +1.36%  "
$ws.Range('E5').Style = "Normal"
$ws.Range('E6').Value = "'  +0.21%  "
$ws.Range('E6').Style = "Normal"
$ws.Range('E7').Value = "'  +1.32%  "
$ws.Range('E7').Style = "Normal"
$ws.Range('D8').Value = "'0.3993"
$ws.Range('D8').Style = "Normal"
$ws.Range('E8').Value = "'  +1.05%  "
$ws.Range('E8').Style = "Normal"
$ws.Range('D9').Value = "'1.522"
$ws.Range('D9').Style = "Normal"
$ws.Range('E9').Value = "'  +4.43%  "
$ws.Range('E9').Style = "Normal"
$ws.Range('D10').Value = "'1.002"
$ws.Range('D10').Style = "Normal"
$ws.Range('E10').Value = "'  +0.31%  "
$ws.Range('E10').Style = "Normal"
$ws.Range('D11').Value = "'52.17"
$ws.Range('D11').Style = "Normal"
$ws.Range('E11').Value = "'  +2.83%  "
$ws.Range('E11').Style = "Normal"
$ws.Range('D12').Value = "'0.08717"
$ws.Range('D12').Style = "Normal"
$ws.Range('E12').Value = "'  +0.65%  "
$ws.Range('E12').Style = "Normal"
$ws.Range('D13').Value = "'7.188"
$ws.Range('D13').Style = "Normal"
$ws.Range('E13').Value = "'  +6.46%  "
$ws.Range('E13').Style = "Normal"
$ws.Range('D14').Value = "'23.01"
$ws.Range('D14').Style = "Normal"
$ws.Range('E14').Value = "'  +1.68%  "
$ws.Range('E14').Style = "Normal"
$ws.Range('D15').Value = "'0.00001313"
$ws.Range('D15').Style = "Normal"
$ws.Range('E15').Value = "'  -0.25%  "
$ws.Range('E15').Style = "Normal"
$ws.Range('D16').Value = "'7.561"
$ws.Range('D16').Style = "Normal"
$ws.Range('E16').Value = "'  +3.80%  "
$ws.Range('E16').Style = "Normal"
$ws.Range('D17').Value = "'1.692.65"
$ws.Range('D17').Style = "Normal"
$ws.Range('E17').Value = "'  +1.59%  "
$ws.Range('E17').Style = "Normal"
$ws.Range('D18').Value = "'99.48"
$ws.Range('D18').Style = "Normal"
$ws.Range('E18').Value = "'  -0.15%  "
$ws.Range('E18').Style = "Normal"
$ws.Range('D19').Value = "'0.07040"
$ws.Range('D19').Style = "Normal"
$ws.Range('E19').Value = "'  +3.77%  "
$ws.Range('E19').Style = "Normal"
$ws.Range('D20').Value = "'19.55"
$ws.Range('D20').Style = "Normal"
$ws.Range('E20').Value = "'  +2.15%  "
$ws.Range('E20').Style = "Normal"
$ws.Range('D21').Value = "'6.855"
$ws.Range('D21').Style = "Normal"
$ws.Range('E21').Value = "'  +3.06%  "
$ws.Range('E21').Style = "Normal"
$ws.Range('E22').Value = "'  +0.11%  "
$ws.Range('E22').Style = "Normal"
$ws.Range('E23').Value = "'  +1.23%  "
$ws.Range('E23').Style = "Normal"
$ws.Range('D24').Value = "'24.517.19"
$ws.Range('D24').Style = "Normal"
$ws.Range('E24').Value = "'  +3.17%  "
$ws.Range('E24').Style = "Normal"
$ws.Range('D25').Value = "'3.066"
$ws.Range('D25').Style = "Normal"
$ws.Range('E25').Value = "'  +7.81%  "
$ws.Range('E25').Style = "Normal"
$ws.Range('D26').Value = "'2.328"
$ws.Range('D26').Style = "Normal"
$ws.Range('E26').Value = "'  +0.61%  "
$ws.Range('E26').Style = "Normal"
$ws.Range('D27').Value = "'22.26"
$ws.Range('D27').Style = "Normal"
$ws.Range('E27').Value = "'  +2.31%  "
$ws.Range('E27').Style = "Normal"
$ws.Range('D28').Value = "'160.81"
$ws.Range('D28').Style = "Normal"
$ws.Range('E28').Value = "'  +0.66%  "
$ws.Range('E28').Style = "Normal"
$ws.Range('D29').Value = "'5.232"
$ws.Range('D29').Style = "Normal"
$ws.Range('E29').Value = "'  +1.22%  "
$ws.Range('E29').Style = "Normal"
$ws.Range('D30').Value = "'133.81"
$ws.Range('D30').Style = "Normal"
$ws.Range('E30').Value = "'  +3.06%  "
$ws.Range('E30').Style = "Normal"
$ws.Range('D31').Value = "'7.495"
$ws.Range('D31').Style = "Normal"
$ws.Range('E31').Value = "'  +10.20%  "
$ws.Range('E31').Style = "Normal"
$ws.Range('D32').Value = "'1.876.42"
$ws.Range('D32').Style = "Normal"
$ws.Range('E32').Value = "'  +1.39%  "
$ws.Range('E32').Style = "Normal"
$ws.Range('D33').Value = "'1.083"
$ws.Range('D33').Style = "Normal"
$ws.Range('E33').Value = "'  -3.36%  "
$ws.Range('E33').Style = "Normal"
$ws.Range('D34').Value = "'0.08517"
$ws.Range('D34').Style = "Normal"
$ws.Range('E34').Value = "'  +0.11%  "
$ws.Range('E34').Style = "Normal"
$ws.Range('D35').Value = "'7.246"
$ws.Range('D35').Style = "Normal"
$ws.Range('E35').Value = "'  +9.77%  "
$ws.Range('E35').Style = "Normal"
$ws.Range('D36').Value = "'11.28"
$ws.Range('D36').Style = "Normal"
$ws.Range('E36').Value = "'  +8.11%  "
$ws.Range('E36').Style = "Normal"
$ws.Range('D37').Value = "'1.945"
$ws.Range('D37').Style = "Normal"
$ws.Range('E37').Value = "'  -0.34%  "
$ws.Range('E37').Style = "Normal"
$ws.Range('D38').Value = "'0.2699"
$ws.Range('D38').Style = "Normal"
$ws.Range('E38').Value = "'  +1.09%  "
$ws.Range('E38').Style = "Normal"
$ws.Range('D39').Value = "'14.38"
$ws.Range('D39').Style = "Normal"
$ws.Range('E39').Value = "'  -0.65%  "
$ws.Range('E39').Style = "Normal"
$ws.Range('D40').Value = "'0.02736"
$ws.Range('D40').Style = "Normal"
$ws.Range('E40').Value = "'  +8.80%  "
$ws.Range('E40').Style = "Normal"
$ws.Range('D41').Value = "'0.09003"
$ws.Range('D41').Style = "Normal"
$ws.Range('E41').Value = "'  +2.34%  "
$ws.Range('E41').Style = "Normal"
$ws.Range('D42').Value = "'1.470"
$ws.Range('D42').Style = "Normal"
$ws.Range('E42').Value = "'  +0.68%  "
$ws.Range('E42').Style = "Normal"
$ws.Range('D43').Value = "'0.7623"
$ws.Range('D43').Style = "Normal"
$ws.Range('E43').Value = "'  +0.57%  "
$ws.Range('E43').Style = "Normal"
$ws.Range('D44').Value = "'0.7140"
$ws.Range('D44').Style = "Normal"
$ws.Range('E44').Value = "'  +1.37%  "
$ws.Range('E44').Style = "Normal"
$ws.Range('D45').Value = "'15.34"
$ws.Range('D45').Style = "Normal"
$ws.Range('E45').Value = "'  +2.40%  "
$ws.Range('E45').Style = "Normal"
$ws.Range('D46').Value = "'2.514"
$ws.Range('D46').Style = "Normal"
$ws.Range('E46').Value = "'  +3.62%  "
$ws.Range('E46').Style = "Normal"
$ws.Range('D47').Value = "'4.194"
$ws.Range('D47').Style = "Normal"
$ws.Range('E47').Value = "'  +2.05%  "
$ws.Range('E47').Style = "Normal"
$ws.Range('E48').Value = "'  +0.18%  "
$ws.Range('E48').Style = "Normal"
$ws.Range('D49').Value = "'140.61"
$ws.Range('D49').Style = "Normal"
$ws.Range('E49').Value = "'  +1.10%  "
$ws.Range('E49').Style = "Normal"
$ws.Range('D50').Value = "'1.326"
$ws.Range('D50').Style = "Normal"
$ws.Range('E50').Value = "'  +7.84%  "
$ws.Range('E50').Style = "Normal"
$ws.Range('D51').Value = "'0.07980"
$ws.Range('D51').Style = "Normal"
$ws.Range('E51').Value = "'  +2.42%  "
$ws.Range('E51').Style = "Normal"
